# 19/09/23 - push 1
#
# This script reproduces (via the Excel COM object model) the edits made to
# "Pasta2.xlsx" between the previous push and this one:
#
#   1. On sheet "exercício 3", the note in A9 (a textual description of the
#      blood-donation eligibility formula) was rewritten.
#   2. On the same sheet, the eligibility formula in column E (E3, and the
#      shared formula starting at E4) was fixed/rewritten so that it
#      correctly checks gender-specific age ranges and the days-since-last
#      donation, returning "apto"/"apta"/"inapto(a)" instead of the old,
#      broken "apto"/"inapto" formula.
#   3. The selection/active cell moved on both "exercício 3" (to J9) and
#      "exercício 4" (to R15), with "exercício 4" ending up as the
#      frontmost (active) sheet/tab.

$wb = $excel.ActiveWorkbook

$wsEx3 = $wb.Worksheets.Item("exercício 3")
$wsEx4 = $wb.Worksheets.Item("exercício 4")

# 1. Update the descriptive note in A9 on "exercício 3"
$wsEx3.Range("A9").Value = "se e(genero, idade, dias); apto; se(e(genero, idade, dias);apto; inapto)"

# 2. Fix the eligibility formula (E3 standalone, E4:E7 shared)
$wsEx3.Range("E3").Formula = '=IF(AND(C3="M",B3>17,B3<66,D3>59),"apto",IF(AND(C3="F",B3>14,B3<61,D3>59),"apta","inapto(a)"))'
$wsEx3.Range("E4:E7").Formula = '=IF(AND(C4="M",B4>17,B4<66,D4>59),"apto",IF(AND(C4="F",B4>14,B4<61,D4>59),"apta","inapto(a)"))'

# 3. Update selections: "exercício 3" -> J9, then "exercício 4" -> R15
#    (selecting "exercício 4" last makes it the active/front sheet tab)
$wsEx3.Select()
$wsEx3.Range("J9").Select()

$wsEx4.Select()
$wsEx4.Range("R15").Select()
